# Update the forecasted "NET" energy values in column B (rows 2-25)
# to reflect the revised MORE Trading Node forecast.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    67215.43566852933,
    63724.96605913646,
    61062.15186725625,
    58690.39048434388,
    57953.41081697737,
    58513.92399294129,
    58560.18761346512,
    61749.84626794878,
    70574.08841398767,
    78215.70537354975,
    82056.68372650631,
    83494.33321479897,
    83824.15966913324,
    87684.05903897769,
    88802.73964341768,
    86996.52656865345,
    84935.95057665437,
    81144.72874415584,
    83246.884137196,
    80548.97919203494,
    78069.16904515574,
    76168.53320712211,
    72298.03743205291,
    67378.63761467961
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $newValues[$i]
}
